$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits right
#    after the title paragraph (it is being relocated/rewritten further down
#    in the document, see step 2).
# ---------------------------------------------------------------------------
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("Meta description")) {
        $metaPara = $para
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Find the trailing "image prompt" paragraph (the italic run that used to
#    read "Create a feature image for Cleopatra ...") and:
#      a) insert a new bold paragraph right before it containing the title
#         text "Play Cleopatra Online Slot Game for Free - Review"
#      b) rewrite its own text to the real meta-description copy, keeping
#         the italic run formatting intact.
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("Create a feature image for Cleopatra")) {
        $targetPara = $para
        break
    }
}

if ($targetPara -ne $null) {
    # 2a. Insert a clean new paragraph ("<empty run/><bold run/>") directly
    #     before the target paragraph. InsertXML's last <w:p> in the package
    #     merges its runs into the paragraph sitting at the insertion point,
    #     so we give it a small text marker run that we strip out right
    #     afterwards -- this avoids leaving a stray duplicate empty run
    #     behind on the untouched paragraph.
    $insertionPoint = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)
    $newParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cleopatra Online Slot Game for Free - Review</w:t></w:r></w:p><w:p><w:r><w:t>@@MARKER@@</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertionPoint.InsertXML($newParaXml)

    $d.Content.Find.Execute("@@MARKER@@", $false, $false, $false, $false, $false, `
                             $true, 1, $false, "", 2) | Out-Null

    # 2b. Re-locate the (still untouched) target paragraph and swap its text
    #     for the new meta-description copy, preserving the italic run.
    $targetPara = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs($i)
        if ($para.Range.Text.StartsWith("Create a feature image for Cleopatra")) {
            $targetPara = $para
            break
        }
    }

    $fullRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End)
    $replacementXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Read our in-depth review of Cleopatra, a popular slot machine game with an Egyptian theme. Play for free on desktop or mobile devices.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $fullRange.InsertXML($replacementXml)
}
